# "Generate Report for Handback"
# d60fcd1c-ad5c-4859-ae55-59a59c812f66 has now been handed back (in both zh-cn and de-de),
# while cf6ce968-bbd2-4664-ad2a-3218565f3e78 remains pending ("Ready for handoff" / "Include").

$wb = $excel.ActiveWorkbook

$dId  = "d60fcd1c-ad5c-4859-ae55-59a59c812f66"
$cId  = "cf6ce968-bbd2-4664-ad2a-3218565f3e78"
$dMd  = "$dId.md"
$cMd  = "$cId.md"
$dZh  = "$dId.efc061ea5a012367f66426a37431b07940d61a60.zh-cn.xlf"
$cZh  = "$cId.1ab83d6b10d671f8826269cb7d7f5f2fd0788f70.zh-cn.xlf"
$dDe  = "$dId.efc061ea5a012367f66426a37431b07940d61a60.de-de.xlf"
$cDe  = "$cId.1ab83d6b10d671f8826269cb7d7f5f2fd0788f70.de-de.xlf"

$dMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/59ab51420053bccfa02e39953c34cd561aedfe32/e2e/$dMd"
$cMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/46b442b4c7f8289d335235cf4dc2cc7420699bae/e2e/$cMd"
$dZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/febd98c3e95120a1f56597b603584d5ab66302b5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$dZh"
$cZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5af2e2e6b124cac4fecb3d77294d6dbea7f10335/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$cZh"
$dDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/26b31d3687b7e87c7b601b6dfa8964dddaeec21b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$dDe"
$cDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/52fbf70eea2eebbc26d9d319f762103ee7823f1e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$cDe"

# ---------------------------------------------------------------------------
# Overview sheet: row order flips so d60fcd1c (now handed back) is listed
# first (row 2) and cf6ce968 (still pending) is listed second (row 3).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-21 14:36:56"

$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-21 14:36:31"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $dMdUrl, [Type]::Missing, [Type]::Missing, $dMd) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $cMdUrl, [Type]::Missing, [Type]::Missing, $cMd) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet: d60fcd1c (row 2) has now been handed back - fill in the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns and flip its Status to "Handed back: in sync with en-US".
# cf6ce968 (row 3) is untouched - still pending handoff.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()

$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("E2").Value = "2016-03-21 14:36:52"
$wsZh.Range("H2").Value = "2016-03-21 14:37:13"
$wsZh.Range("J2").Value = "Include"

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = "2016-03-21 14:36:28"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("J3").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $dMdUrl, [Type]::Missing, [Type]::Missing, $dMd) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $dZhUrl, [Type]::Missing, [Type]::Missing, $dZh) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $dMdUrl, [Type]::Missing, [Type]::Missing, $dMd) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $dZhUrl, [Type]::Missing, [Type]::Missing, $dZh) | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $cMdUrl, [Type]::Missing, [Type]::Missing, $cMd) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $cZhUrl, [Type]::Missing, [Type]::Missing, $cZh) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: same story as zh-cn but for the de-de locale.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()

$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("E2").Value = "2016-03-21 14:36:56"
$wsDe.Range("H2").Value = "2016-03-21 14:37:22"
$wsDe.Range("J2").Value = "Include"

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("E3").Value = "2016-03-21 14:36:31"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("J3").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $dMdUrl, [Type]::Missing, [Type]::Missing, $dMd) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $dDeUrl, [Type]::Missing, [Type]::Missing, $dDe) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $dMdUrl, [Type]::Missing, [Type]::Missing, $dMd) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $dDeUrl, [Type]::Missing, [Type]::Missing, $dDe) | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $cMdUrl, [Type]::Missing, [Type]::Missing, $cMd) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $cDeUrl, [Type]::Missing, [Type]::Missing, $cDe) | Out-Null

Write-Host "Report generated for handback."
